$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New helper data in column F (used by the SUBTOTAL formula below) ---
$ws.Range("F1").Value = 1
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 1

# --- Labels describing the two new conditional-formatting demo rows ---
$ws.Range("A3").Value = "cell with formula  condition"
$ws.Range("A4").Value = "cell with formula value and > condition"

# B3 previously triggered a NotImplementedException because of the SUBTOTAL
# formula referenced by a conditional-formatting expression rule (SHEET-307).
$ws.Range("B3").Formula = "=SUBTOTAL(109,F1:F5)"
$ws.Range("B4").Value = 5

# --- Make column A wide enough to show the labels ---
$ws.Columns("A").ColumnWidth = 31.25

# --- Conditional formatting: B4 highlights when less than B3 ---
$cfB4 = $ws.Range("B4").FormatConditions.Add(1, 6, '=$B$3')
$cfB4.Font.Color = -16383844
$cfB4.Interior.Color = 13551615

# --- Conditional formatting: B3 highlights when the expression $B$3>5 is true ---
$cfB3 = $ws.Range("B3").FormatConditions.Add(2, 0, '=$B$3>5')
$cfB3.Font.Color = -16383844
$cfB3.Interior.Color = 13551615

# --- Move the active selection to B4, matching the saved view state ---
$ws.Range("B4").Select()

Write-Host "done"
